# Add a new "Turkey" country-template sheet, cloned from the existing
# "Spain" sheet (same layout/styles), trimmed down to a single data row,
# and re-pointed at a new "NGC-3191/T3329" identifier string.

$wb = $excel.ActiveWorkbook

$spain = $wb.Worksheets.Item("Spain")

# Leave Spain's own selection the way the author left it before copying.
$null = $spain.Range("A1:W9").Select()

# Move-or-Copy "Spain" to a new sheet placed right after it (same as
# right-click > Move or Copy... > Create a copy, dropped after "Spain").
$null = $spain.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

$turkey = $wb.Worksheets.Item($wb.Worksheets.Count)
$turkey.Name = "Turkey"

# Only one loading-detail data row is needed for the new template, so
# drop the extra two example rows (9 and 10) that came over from Spain.
$turkey.Rows.Item(9).EntireRow.Delete()
$turkey.Rows.Item(9).EntireRow.Delete()

# Point the template reference cell at the new Turkey identifier.
$turkey.Range("B4").Value = "NGC-3191/T3329"

# Re-fit the row heights for the wrapped-text rows (these end up
# auto-sized in Excel once the extra rows are gone / content changes).
$turkey.Rows.Item(3).RowHeight = 28.8
$turkey.Rows.Item(4).RowHeight = 28.8
$turkey.Rows.Item(5).RowHeight = 28.8
$turkey.Rows.Item(7).RowHeight = 72

# Columns A and B best-fit to their (now narrower) content instead of
# the wide layout columns that were specific to Spain's sheet.
$turkey.Columns.Item(1).ColumnWidth = 24.3
$turkey.Columns.Item(2).ColumnWidth = 14.4
$turkey.Columns.Item(12).ColumnWidth = 14
$turkey.Columns.Item(14).ColumnWidth = 14
$turkey.Columns.Item(16).ColumnWidth = 14
$turkey.Columns.Item(18).ColumnWidth = 14

# Leave the new sheet selected near the data table, like the author did.
$null = $turkey.Range("H11").Select()

Write-Host "Added Turkey sheet"
